# Daily attendance processing - 2025-12-03 01:27:45
# Swap the order of the first two comma-separated "Recorded By" entries
# in column G for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Text

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -ge 2) {
            $first = $parts[0]
            $second = $parts[1]
            $parts[0] = $second
            $parts[1] = $first
            $newVal = [string]::Join(", ", $parts)
            $cell.Value = $newVal
        }
    }
}
